# Added HomePage Test Suites
# Adds three new test-case rows (verifyChooseModemTwo/Three/Four) to the
# "Testdata" worksheet, mirroring the formatting of the existing
# "verifyChooseModemOne" row (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testdata")

# New test case data: TCName, text1, text2, text3, text4(nickname)
$newRows = @(
    @{ Row = 6; TCName = "verifyChooseModemTwo";   Text1 = "text1=Welcometo Globe myBusiness"; Text2 = "mobileNumber=09271080510"; Text3 = "pin=1111"; Text4 = "nickname=Abraham" },
    @{ Row = 7; TCName = "verifyChooseModemThree"; Text1 = "text1=Welcometo Globe myBusiness"; Text2 = "mobileNumber=09271080510"; Text3 = "pin=1111"; Text4 = "nickname=Ham" },
    @{ Row = 8; TCName = "verifyChooseModemFour";  Text1 = "text1=Welcometo Globe myBusiness"; Text2 = "mobileNumber=09271080510"; Text3 = "pin=1111"; Text4 = "nickname=Honestabe" }
)

foreach ($rowInfo in $newRows) {
    $r = $rowInfo.Row

    # Copy the formatting from the reference row (row 5, columns A:E) so the
    # new rows look exactly like the existing "verifyChooseModemOne" row.
    $ws.Range("A5:E5").Copy() | Out-Null
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $rowInfo.TCName
    $ws.Cells.Item($r, 2).Value = $rowInfo.Text1
    $ws.Cells.Item($r, 3).Value = $rowInfo.Text2
    $ws.Cells.Item($r, 4).Value = $rowInfo.Text3
    $ws.Cells.Item($r, 5).Value = $rowInfo.Text4
}

$excel.CutCopyMode = 0
